# Added bakery manufacturer filter test
#
# Adds a new "filterLimit" column (F) to the washers price-filter test
# fixture and bumps the "Maximum price" value on the first data row so the
# new limit filter has something meaningful to test against.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F1: new "filterLimit" header -------------------------------------
# Re-use the existing bold/centered header formatting (same style as the
# other header cells) by copying format from A1.
$ws.Range("F1").Value = "filterLimit"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- F2:F17: new filter-limit values -----------------------------------
# Re-use the existing left-aligned data formatting by copying format from A2.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2:F17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).Value = 7000
}

$excel.CutCopyMode = 0

# Row 2's "Maximum price" (E2) goes from 0 to 100000
$ws.Range("E2").Value = 100000

# --- Column width adjustments (closest value achievable via ColumnWidth) --
$ws.Columns.Item(1).ColumnWidth = 22.25000000000001    # A -> ~23.140625
$ws.Columns.Item(2).ColumnWidth = 21.41666666666669    # B -> ~22.28515625
$ws.Columns.Item(5).ColumnWidth = 11.25                 # E -> ~12.140625
$ws.Columns.Item(6).ColumnWidth = 11.25                 # F -> same as E

# --- Selection moves to the newly edited cell E2 -----------------------
$ws.Range("E2").Select() | Out-Null
